# [맘파] RecipeProbabilityTable ID누락 수정
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RecipeProbability")

# Column B width fix (target OOXML width 14.14; the COM ColumnWidth property
# is character-unit and gets snapped to the host's pixel grid on save, so
# 13.3 is the closest input that round-trips to the nearest achievable width)
$ws.Columns.Item(2).ColumnWidth = 13.3

# Probability / GroupID corrections
$ws.Range("C6").Value = 19.0
$ws.Range("C7").Value = 1.0
$ws.Range("C9").Value = 19.0
$ws.Range("D10").Value = 20.0
$ws.Range("C11").Value = 1.0
$ws.Range("D11").Value = 20.0
$ws.Range("C12").Value = 19.0
$ws.Range("D12").Value = 20.0
$ws.Range("C13").Value = 5.0
$ws.Range("D13").Value = 20.0
$ws.Range("C14").Value = 7.0
$ws.Range("D14").Value = 20.0
